# deployment and gantt chart formatting
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "Baseline-Mid (SC)"

# New Year / Cumulative Capacity data (rows 2-17), replacing the old rows 2-27
$years = @(2031,2032,2033,2034,2035,2036,2037,2038,2039,2040,2041,2042,2043,2044,2045,2046)
$values = @(
    1014.686248331108,
    1723.928200563715,
    2652.480900977594,
    3737.659267944707,
    5405.037494671253,
    7286.013900313775,
    9144.594764214053,
    10954.11508666184,
    12768.9411359356,
    14623.6582912613,
    16559.53882147466,
    18470.18061703362,
    20290.75680477099,
    22268.49720168707,
    23376.74780385629,
    24496.32143330742
)

for ($i = 0; $i -lt $years.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $years[$i]
    $ws.Cells.Item($row, 2).Value = $values[$i]
}

# Remove the now-unused trailing rows (old rows 18-27, i.e. years 2047-2057)
$oldLastRow = 27
$newLastRow = 17
if ($oldLastRow -gt $newLastRow) {
    $clearRange = $ws.Range($ws.Cells.Item($newLastRow + 1, 1), $ws.Cells.Item($oldLastRow, 2))
    $clearRange.Clear()
}
